$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.611.37'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.850.96'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '264.43'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5219'
$ws.Range('E7').Value = '  +1.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3267'
$ws.Range('E8').Value = '  +0.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06782'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.80'
$ws.Range('E10').Value = '  -1.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7750'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07757'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.854.29'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.41'
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.028'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.99'
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007943'
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.654.88'
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.094.05'
$ws.Range('E21').Value = '  -1.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.639'
$ws.Range('E22').Value = '  +1.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.530'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.002'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.64'
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.200'
$ws.Range('E26').Value = '  -5.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.675'
$ws.Range('E27').Value = '  +1.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.04'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.09'
$ws.Range('E29').Value = '  +1.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.201'
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.137'
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08763'
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04829'
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.134'
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.876'
$ws.Range('E35').Value = '  +0.53%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7148'
$ws.Range('E36').Value = '  +5.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.110'
$ws.Range('E37').Value = '  +1.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01790'
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.199'
$ws.Range('E39').Value = '  -1.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4877'
$ws.Range('E40').Value = '  -1.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '112.71'
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8953'
$ws.Range('E42').Value = '  -0.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.075'
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.704'
$ws.Range('E45').Value = '  -0.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4179'
$ws.Range('E46').Value = '  -1.76%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.128'
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05922'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.05'
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1232'
$ws.Range('E50').Value = '  -4.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8852'
$ws.Range('E51').Value = '  +3.51%  '
